# Update filtered_output.xlsx: refresh Neg_Change (grows from 4 to 13 data rows)
# and Pos_Change (stays at 5 data rows, values refreshed) with new market data.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Neg_Change sheet
# ---------------------------------------------------------------------------
$wsNeg = $wb.Worksheets.Item("Neg_Change")

$negRows = @(
    @("ICICIBANK", 1372.3, 1381.6, 1364.3, 1367.8, 8100832, 16126464, -0.497668428739245, "ICICIBANK"),
    @("COALINDIA", 385.5, 386.4, 381.1, 382.05, 3642904, 7895734, -0.5386237682272478, "COALINDIA"),
    @("JSWSTEEL", 1158, 1162.7, 1142.7, 1146, 737972, 1565661, -0.5286514769161396, "JSWSTEEL"),
    @("ONGC", 246.55, 246.8, 241.14, 241.99, 5073558, 12498420, -0.5940640496958816, "ONGC"),
    @("JIOFIN", 309.9, 310.95, 304.2, 304.7, 7853613, 18717966, -0.5804238024580235, "JIOFIN"),
    @("ADANIENSOL", 923.45, 926.9, 908.15, 913.3, 662712, 1341521, -0.5059995333654859, "ADANIENSOL"),
    @("HINDPETRO", 452, 463, 452, 458, 4192443, 8444885, -0.50355238703665, "HINDPETRO"),
    @("PATANJALI", 596.95, 599.9, 593.1, 596, 185683, 413182, -0.5506023979747424, "PATANJALI"),
    @("BANKINDIA", 126.65, 126.69, 123.7, 125, 7459860, 17946269, -0.5843225129412692, "BANKINDIA"),
    @("ATGL", 638, 638.4, 627.2, 627.2, 742143, 1784826, -0.5841930810062157, "ATGL"),
    @("AUROPHARMA", 1091.5, 1102.8, 1066.8, 1068, 549749, 1177833, -0.5332538653612184, "AUROPHARMA"),
    @("DELHIVERY", 470.2, 471.5, 461.05, 468.3, 2249868, 4836942, -0.5348573540885957, "DELHIVERY"),
    @("ANGELONE", 2242, 2259, 2210.2, 2213, 555881, 1137247, -0.511204689922242, "ANGELONE")
)

$r = 2
foreach ($row in $negRows) {
    $c = 1
    foreach ($val in $row) {
        $wsNeg.Cells.Item($r, $c).Value = $val
        $c = $c + 1
    }
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# Pos_Change sheet
# ---------------------------------------------------------------------------
$wsPos = $wb.Worksheets.Item("Pos_Change")

$posRows = @(
    @("SBILIFE", 1784, 1795.9, 1757.7, 1772, 1323877, 859968, 0.5394491422936667, "SBILIFE"),
    @("ASIANPAINT", 2360, 2369.4, 2320, 2323, 701799, 498524, 0.4077536888896021, "ASIANPAINT"),
    @("PIDILITIND", 1469, 1498.5, 1462, 1495, 1023522, 725776, 0.4102450342805494, "PIDILITIND"),
    @("COFORGE", 1682.5, 1733.2, 1673.3, 1709.5, 2733508, 1937983, 0.4104912168992195, "COFORGE"),
    @("KEI", 4175, 4297, 4170.2, 4250, 225852, 143928, 0.5692012673003168, "KEI")
)

$r = 2
foreach ($row in $posRows) {
    $c = 1
    foreach ($val in $row) {
        $wsPos.Cells.Item($r, $c).Value = $val
        $c = $c + 1
    }
    $r = $r + 1
}
